$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency price/volume snapshot values.
# A leading apostrophe forces Excel to store numeric-looking strings
# (e.g. "574.96", "25.20") as text instead of coercing them to numbers,
# which preserves exact formatting such as trailing zeros.

$ws.Range("D2").Value = '61.301.09'
$ws.Range("E2").Value = '  -2.30%  '
$ws.Range("D3").Value = '3.384.24'
$ws.Range("E3").Value = '  -0.17%  '
$ws.Range("D4").Value = '''0.999'
$ws.Range("E4").Value = '  -0.07%  '
$ws.Range("D5").Value = '''574.96'
$ws.Range("E5").Value = '  +0.22%  '
$ws.Range("D6").Value = '''135.72'
$ws.Range("E6").Value = '  +7.70%  '
$ws.Range("E7").Value = '  -0.06%  '
$ws.Range("D8").Value = '3.383.31'
$ws.Range("E8").Value = '  -0.20%  '
$ws.Range("E9").Value = '  +0.89%  '
$ws.Range("E10").Value = '  +3.20%  '
$ws.Range("E11").Value = '  +1.90%  '
$ws.Range("D12").Value = '''0.389'
$ws.Range("E12").Value = '  +2.78%  '
$ws.Range("D13").Value = '3.958.84'
$ws.Range("E13").Value = '  -0.22%  '
$ws.Range("E14").Value = '  +1.07%  '
$ws.Range("E15").Value = '  +0.99%  '
$ws.Range("D16").Value = '3.383.57'
$ws.Range("E16").Value = '  -0.16%  '
$ws.Range("D17").Value = '''25.20'
$ws.Range("E17").Value = '  +1.66%  '
$ws.Range("D18").Value = '61.382.03'
$ws.Range("E18").Value = '  -2.17%  '
$ws.Range("D19").Value = '''14.08'
$ws.Range("E19").Value = '  +6.96%  '
$ws.Range("D20").Value = '''5.81'
$ws.Range("E20").Value = '  +2.49%  '
$ws.Range("D21").Value = '''9.41'
$ws.Range("E21").Value = '  -1.07%  '
$ws.Range("D22").Value = '''376.84'
$ws.Range("E22").Value = '  +0.22%  '
$ws.Range("D23").Value = '''0.569'
$ws.Range("E23").Value = '  +2.07%  '
$ws.Range("D24").Value = '3.514.04'
$ws.Range("E24").Value = '  -0.39%  '
$ws.Range("E25").Value = '  +0.05%  '
$ws.Range("D26").Value = '''70.68'
$ws.Range("E26").Value = '  -2.09%  '
$ws.Range("E27").Value = '  +9.87%  '
$ws.Range("D28").Value = '''1.70'
$ws.Range("E28").Value = '  +22.00%  '
$ws.Range("D29").Value = '''7.82'
$ws.Range("E29").Value = '  +12.22%  '
$ws.Range("E30").Value = '  +0.12%  '
$ws.Range("D31").Value = '''8.16'
$ws.Range("E31").Value = '  +4.16%  '
$ws.Range("B32").Value = 'PancakeSwap'
$ws.Range("C32").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D32").Value = '''2.17'
$ws.Range("E32").Value = '  +0.63%  '
$ws.Range("B33").Value = 'Kaspa'
$ws.Range("C33").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D33").Value = '''0.156'
$ws.Range("E33").Value = '  +4.15%  '
$ws.Range("E34").Value = '  -0.08%  '
$ws.Range("D35").Value = '3.412.23'
$ws.Range("E35").Value = '  -0.23%  '
$ws.Range("D36").Value = '''23.45'
$ws.Range("E36").Value = '  +3.18%  '
$ws.Range("D37").Value = '''5.61'
$ws.Range("E37").Value = '  +6.47%  '
$ws.Range("D38").Value = '''6.98'
$ws.Range("E38").Value = '  +3.88%  '
$ws.Range("D39").Value = '''1.57'
$ws.Range("E39").Value = '  +5.53%  '
$ws.Range("D40").Value = '''162.13'
$ws.Range("D41").Value = '''0.0788'
$ws.Range("E41").Value = '  +4.03%  '
$ws.Range("E42").Value = '  -0.01%  '
$ws.Range("E43").Value = '  +13.04%  '
$ws.Range("D44").Value = '''4.44'
$ws.Range("E44").Value = '  +4.01%  '
$ws.Range("D45").Value = '''41.64'
$ws.Range("E45").Value = '  +0.41%  '
$ws.Range("D46").Value = '''0.762'
$ws.Range("E46").Value = '  -1.55%  '
$ws.Range("D47").Value = '''1.62'
$ws.Range("E47").Value = '  +3.00%  '
$ws.Range("D48").Value = '''23.72'
$ws.Range("E48").Value = '  +3.47%  '
$ws.Range("D49").Value = '''6.96'
$ws.Range("E49").Value = '  +4.86%  '
$ws.Range("D50").Value = '''23.06'
$ws.Range("E50").Value = '  +13.72%  '
$ws.Range("D51").Value = '''0.904'
$ws.Range("E51").Value = '  +6.24%  '
